$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 83336880
$ws.Range("I62").Value = 2342.7144
$ws.Range("J62").Value = 200005250
$ws.Range("K62").Value = 2342.7144
$ws.Range("L62").Value = 200005250
$ws.Range("M62").Value = -1718.7144
$ws.Range("N62").Value = -200006498

$ws.Range("H65").Value = 83336880
$ws.Range("I65").Value = 2342.7144
$ws.Range("J65").Value = 200005250
$ws.Range("K65").Value = 11713.572
$ws.Range("L65").Value = 1000026250
$ws.Range("M65").Value = -8593.572
$ws.Range("N65").Value = -1000032490

$ws.Range("H113").Value = 2522.0857
$ws.Range("I113").Value = 2002.5
$ws.Range("J113").Value = 2553.5757
$ws.Range("K113").Value = 2002.5
$ws.Range("L113").Value = 2553.5757
$ws.Range("M113").Value = 1251.5
$ws.Range("N113").Value = -9061.575699999999

$ws.Range("H117").Value = 56000
$ws.Range("J117").Value = 56000
$ws.Range("L117").Value = 56000
$ws.Range("N117").Value = -65178

$ws.Range("H129").Value = 1812.9025
$ws.Range("I129").Value = 336.125
$ws.Range("J129").Value = 2170.9092
$ws.Range("K129").Value = 1008.375
$ws.Range("L129").Value = 6512.7276
$ws.Range("M129").Value = 3991.625
$ws.Range("N129").Value = -16512.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 100002664
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 111113900
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 111113900
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -111114712

$ws.Range("H91").Value = 100002664
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 111113900
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 111113900
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -111116708

$ws.Range("H132").Value = 2105.9714
$ws.Range("I132").Value = 1822.0358
$ws.Range("J132").Value = 3241.7144
$ws.Range("K132").Value = 5466.107400000001
$ws.Range("L132").Value = 9725.143199999999
$ws.Range("M132").Value = -2936.107400000001
$ws.Range("N132").Value = -14785.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9525769
$ws.Range("I86").Value = 12501985
$ws.Range("J86").Value = 1876
$ws.Range("K86").Value = 12501985
$ws.Range("L86").Value = 1876
$ws.Range("M86").Value = -12500862
$ws.Range("N86").Value = -4122

$ws.Range("H89").Value = 9525769
$ws.Range("I89").Value = 12501985
$ws.Range("J89").Value = 1876
$ws.Range("K89").Value = 62509925
$ws.Range("L89").Value = 9380
$ws.Range("M89").Value = -62504309
$ws.Range("N89").Value = -20612

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6267404
$ws.Range("I31").Value = 4903252.5
$ws.Range("J31").Value = 8699152
$ws.Range("K31").Value = 4903252.5
$ws.Range("L31").Value = 8699152
$ws.Range("M31").Value = -4902957.5
$ws.Range("N31").Value = -8699742

$ws.Range("H34").Value = 6267404
$ws.Range("I34").Value = 4903252.5
$ws.Range("J34").Value = 8699152
$ws.Range("K34").Value = 4903252.5
$ws.Range("L34").Value = 8699152
$ws.Range("M34").Value = -4903050.5
$ws.Range("N34").Value = -8699556

$ws.Range("H62").Value = 83338024
$ws.Range("I62").Value = 3000.8333
$ws.Range("J62").Value = 166673040
$ws.Range("K62").Value = 3000.8333
$ws.Range("L62").Value = 166673040
$ws.Range("M62").Value = -2376.8333
$ws.Range("N62").Value = -166674288

$ws.Range("H65").Value = 83338024
$ws.Range("I65").Value = 3000.8333
$ws.Range("J65").Value = 166673040
$ws.Range("K65").Value = 15004.1665
$ws.Range("L65").Value = 833365200
$ws.Range("M65").Value = -11884.1665
$ws.Range("N65").Value = -833371440

$ws.Range("H99").Value = 1155526.5
$ws.Range("I99").Value = 1704024.2
$ws.Range("J99").Value = 3681.4
$ws.Range("K99").Value = 1704024.2
$ws.Range("L99").Value = 3681.4
$ws.Range("M99").Value = -1702526.2
$ws.Range("N99").Value = -6677.4

$ws.Range("H126").Value = 1155526.5
$ws.Range("I126").Value = 1704024.2
$ws.Range("J126").Value = 3681.4
$ws.Range("K126").Value = 5112072.6
$ws.Range("L126").Value = 11044.2
$ws.Range("M126").Value = -5109602.6
$ws.Range("N126").Value = -15984.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 906.8461
$ws.Range("I5").Value = 880.8182
$ws.Range("J5").Value = 1050
$ws.Range("K5").Value = 2642.4546
$ws.Range("L5").Value = 3150
$ws.Range("M5").Value = -2530.4546
$ws.Range("N5").Value = -3374

$ws.Range("H122").Value = 637.6667
$ws.Range("I122").Value = 413.6
$ws.Range("J122").Value = 1085.8
$ws.Range("K122").Value = 3722.4
$ws.Range("L122").Value = 9772.199999999999
$ws.Range("M122").Value = -1272.4
$ws.Range("N122").Value = -14672.2

$ws.Range("H135").Value = 906.8461
$ws.Range("I135").Value = 880.8182
$ws.Range("J135").Value = 1050
$ws.Range("K135").Value = 7927.3638
$ws.Range("L135").Value = 9450
$ws.Range("M135").Value = -5392.3638
$ws.Range("N135").Value = -14520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2466.4211
$ws.Range("I132").Value = 1726.9231
$ws.Range("J132").Value = 4068.6667
$ws.Range("K132").Value = 5180.7693
$ws.Range("L132").Value = 12206.0001
$ws.Range("M132").Value = -2650.7693
$ws.Range("N132").Value = -17266.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1884.2667
$ws.Range("I7").Value = 1525.375
$ws.Range("J7").Value = 2294.4285
$ws.Range("K7").Value = 1525.375
$ws.Range("L7").Value = 2294.4285
$ws.Range("M7").Value = -1413.375
$ws.Range("N7").Value = -2518.4285

$ws.Range("H40").Value = 1989.2307
$ws.Range("I40").Value = 1914.2858
$ws.Range("J40").Value = 2076.6667
$ws.Range("K40").Value = 1914.2858
$ws.Range("L40").Value = 2076.6667
$ws.Range("M40").Value = -1778.2858
$ws.Range("N40").Value = -2348.6667

$ws.Range("H126").Value = 1884.2667
$ws.Range("I126").Value = 1525.375
$ws.Range("J126").Value = 2294.4285
$ws.Range("K126").Value = 4576.125
$ws.Range("L126").Value = 6883.2855
$ws.Range("M126").Value = -2106.125
$ws.Range("N126").Value = -11823.2855

$ws.Range("H132").Value = 1847988.5
$ws.Range("I132").Value = 3691781
$ws.Range("J132").Value = 4196.0586
$ws.Range("K132").Value = 11075343
$ws.Range("L132").Value = 12588.1758
$ws.Range("M132").Value = -11072813
$ws.Range("N132").Value = -17648.1758
